# Add new chatbot intents: TensorFlow, React, NLP, and a second "tabular data"
# style row, appended after the existing data (rows 1-79) in "Intents" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 80: ask.what_is_tensorflow ---------------------------------------
$ws.Range("A80").Value = "ask.what_is_tensorflow"
$ws.Range("B80").Value = "what is tensorflow?; describe tensorflow?; tell me something about Tensorflow?; define Tensorflow?; explain tensorflow?;"
$ws.Range("C80").Value = "TensorFlow is a free and open-source software library for machine learning and artificial intelligence. It can be used across a range of tasks but has a particular focus on training and inference of deep neural networks. TensorFlow provides a collection of workflows to develop and train models using Python or JavaScript, and to easily deploy in the cloud, on-prem, in the browser, or on-device no matter what language you use. The tf. data API enables you to build complex input pipelines from simple, reusable pieces."

# --- Row 81: ask.what_is_react (name entered ahead of its Q/A text) -------
$ws.Range("A81").Value = "ask.what_is_react"

# --- Row 82: ask.what_is_nlp (name entered ahead of its Q/A text) ---------
$ws.Range("A82").Value = "ask.what_is_nlp"

# --- Row 83: reuse of the existing "tabular data" intent, new answer ------
$ws.Range("C83").Value = "A table is an arrangement of information or data, typically in rows and columns, or possibly in a more complex structure. Tables are widely used in communication, research, and data analysis."

# --- Fill in remaining React / NLP text ------------------------------------
$ws.Range("B81").Value = "what is react?; what is react.js; what is React JS?; do you know react js?;"
$ws.Range("B82").Value = "what is nlp?; describe NLP?; define NLP; explain NLP; can you explain NLP?;"
$ws.Range("C82").Value = "NLP or Natural language processing is a subfield of linguistics, computer science, and artificial intelligence concerned with the interactions between computers and human language, in particular how to program computers to process and analyze large amounts of natural language data."
$ws.Range("C81").Value = "React is a free and open-source front-end JavaScript library for building user interfaces based on UI components. It is maintained by Meta and a community of individual developers and companies. React can be used as a base in the development of single-page or mobile applications."

$ws.Range("A83").Value = "ask.what_is_tabular_data"
$ws.Range("B83").Value = "What is tabular data?;"

# Column B got a bit wider to fit the new (longer) questions.
$ws.Columns("B").ColumnWidth = 113.42578125

# Reflect where the user ended up after typing all this in (also scrolls
# the view so C61 becomes the top-left visible cell).
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C91").Select()
